$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4058.2856
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4058.2856
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H69").Value = 9833.166999999999
$ws.Range("I69").Value = 8800
$ws.Range("K69").Value = 26400
$ws.Range("M69").Value = -25526
$ws.Range("H72").Value = 9833.166999999999
$ws.Range("I72").Value = 8800
$ws.Range("K72").Value = 79200
$ws.Range("M72").Value = -74832
$ws.Range("H86").Value = 14320.1
$ws.Range("I86").Value = 13740.4
$ws.Range("K86").Value = 13740.4
$ws.Range("M86").Value = -12617.4
$ws.Range("H89").Value = 14320.1
$ws.Range("I89").Value = 13740.4
$ws.Range("K89").Value = 68702
$ws.Range("M89").Value = -63086
$ws.Range("H112").Value = 2774.3333
$ws.Range("I112").Value = 2937.8
$ws.Range("J112").Value = 2657.5715
$ws.Range("K112").Value = 8813.400000000001
$ws.Range("L112").Value = 7972.7145
$ws.Range("M112").Value = -7705.400000000001
$ws.Range("N112").Value = -10188.7145
$ws.Range("H131").Value = 1972.9
$ws.Range("I131").Value = 1914.3334
$ws.Range("J131").Value = 2500
$ws.Range("K131").Value = 5743.0002
$ws.Range("L131").Value = 7500
$ws.Range("M131").Value = -703.0002000000004
$ws.Range("N131").Value = -17580
$ws.Range("H132").Value = 27591
$ws.Range("I132").Value = 31246.857
$ws.Range("K132").Value = 93740.571
$ws.Range("M132").Value = -91210.571
$ws.Range("H135").Value = 6092.5557
$ws.Range("I135").Value = 2272.818
$ws.Range("K135").Value = 20455.362
$ws.Range("M135").Value = -17920.362
$ws.Range("H138").Value = 3964.24
$ws.Range("I138").Value = 3129.7334
$ws.Range("K138").Value = 9389.200199999999
$ws.Range("M138").Value = -4249.200199999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 83730.55499999999
$ws.Range("J139").Value = 83730.55499999999
$ws.Range("L139").Value = 83730.55499999999
$ws.Range("N139").Value = -94010.55499999999
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 53662.45
$ws.Range("I20").Value = 70451.664
$ws.Range("K20").Value = 70451.664
$ws.Range("M20").Value = -70204.664
$ws.Range("H22").Value = 343.375
$ws.Range("I22").Value = 360.5
$ws.Range("J22").Value = 292
$ws.Range("K22").Value = 360.5
$ws.Range("L22").Value = 292
$ws.Range("M22").Value = -187.5
$ws.Range("N22").Value = -638
$ws.Range("H26").Value = 23931.834
$ws.Range("I26").Value = 15718.2
$ws.Range("K26").Value = 15718.2
$ws.Range("M26").Value = -15426.2
$ws.Range("H96").Value = 37566.832
$ws.Range("J96").Value = 53333.332
$ws.Range("L96").Value = 53333.332
$ws.Range("N96").Value = -58825.332
$ws.Range("H105").Value = 6439.5
$ws.Range("I105").Value = 1929.6666
$ws.Range("K105").Value = 1929.6666
$ws.Range("M105").Value = -182.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6200
$ws.Range("I31").Value = 6200
$ws.Range("K31").Value = 6200
$ws.Range("M31").Value = -5905
$ws.Range("H34").Value = 6200
$ws.Range("I34").Value = 6200
$ws.Range("K34").Value = 6200
$ws.Range("M34").Value = -5998
$ws.Range("H107").Value = 1225.7037
$ws.Range("I107").Value = 1124.2142
$ws.Range("K107").Value = 1124.2142
$ws.Range("M107").Value = 795.7858000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1750.7646
$ws.Range("J34").Value = 1879.4615
$ws.Range("L34").Value = 5638.3845
$ws.Range("N34").Value = -5806.3845
$ws.Range("H37").Value = 44998.57
$ws.Range("J37").Value = 44998.57
$ws.Range("L37").Value = 134995.71
$ws.Range("N37").Value = -135219.71
$ws.Range("H128").Value = 324769.34
$ws.Range("I128").Value = 324769.34
$ws.Range("K128").Value = 974308.02
$ws.Range("M128").Value = -969328.02
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2508.5
$ws.Range("I43").Value = 2508.5
$ws.Range("K43").Value = 2508.5
$ws.Range("M43").Value = -2357.5
$ws.Range("H46").Value = 12199.5
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 19399
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 19399
$ws.Range("M46").Value = -4844
$ws.Range("N46").Value = -19711
$ws.Range("H57").Value = 20599.334
$ws.Range("J57").Value = 22719.2
$ws.Range("L57").Value = 22719.2
$ws.Range("N57").Value = -24359.2
$ws.Range("H80").Value = 4407.857
$ws.Range("I80").Value = 4701.6665
$ws.Range("K80").Value = 4701.6665
$ws.Range("M80").Value = -3703.6665
$ws.Range("H83").Value = 4407.857
$ws.Range("I83").Value = 4701.6665
$ws.Range("K83").Value = 23508.3325
$ws.Range("M83").Value = -18516.3325
$ws.Range("H122").Value = 2447.8
$ws.Range("I122").Value = 1640
$ws.Range("J122").Value = 2986.3333
$ws.Range("K122").Value = 4920
$ws.Range("L122").Value = 8958.999899999999
$ws.Range("M122").Value = -2470
$ws.Range("N122").Value = -13858.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4339.3335
$ws.Range("I46").Value = 1154.75
$ws.Range("J46").Value = 5249.2144
$ws.Range("K46").Value = 1154.75
$ws.Range("L46").Value = 5249.2144
$ws.Range("M46").Value = -966.75
$ws.Range("N46").Value = -5625.2144
$ws.Range("H82").Value = 2235
$ws.Range("J82").Value = 2266.6667
$ws.Range("L82").Value = 2266.6667
$ws.Range("N82").Value = -2988.6667
$ws.Range("H85").Value = 2235
$ws.Range("J85").Value = 2266.6667
$ws.Range("L85").Value = 2266.6667
$ws.Range("N85").Value = -4762.6667
$ws.Range("H94").Value = 88999.5
$ws.Range("J94").Value = 88999.5
$ws.Range("L94").Value = 88999.5
$ws.Range("N94").Value = -90351.5
$ws.Range("H132").Value = 3662.9
$ws.Range("I132").Value = 2580
$ws.Range("J132").Value = 6189.6665
$ws.Range("K132").Value = 7740
$ws.Range("L132").Value = 18568.9995
$ws.Range("M132").Value = -5210
$ws.Range("N132").Value = -23628.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1242.2046
$ws.Range("J107").Value = 1876.75
$ws.Range("L107").Value = 5630.25
$ws.Range("N107").Value = -9470.25
$ws.Range("H122").Value = 90892.16
$ws.Range("I122").Value = 4108.5557
$ws.Range("J122").Value = 286155.25
$ws.Range("K122").Value = 12325.6671
$ws.Range("L122").Value = 858465.75
$ws.Range("M122").Value = -9875.667099999999
$ws.Range("N122").Value = -863365.75
